$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the document (after "Ошмяны
#    20 23") to wrap the "В. Ф. Одиночко" signature text, and make that
#    signature text bold (mirrors the author re-visiting / re-typing that
#    run last, which is what moves Word's internal "last edit" bookmark).
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$found = $rng.Find.Execute("В. Ф. Одиночко")
if ($found) {
  $rng.Bold = 1
  $d.Bookmarks.Add("_GoBack", $rng)
}

# ---------------------------------------------------------------------------
# 2) Register the "Balloon Text" paragraph style (a6) and its linked
#    character style "Текст выноски Знак" (a7) in styles.xml, matching the
#    pair Word mints automatically. We mint them through a throwaway
#    paragraph at the very end of the story so the real content is
#    untouched, then immediately remove that scratch paragraph again.
# ---------------------------------------------------------------------------

$scratchRange = $d.Content
$scratchRange.Collapse(0)
$scratchRange.InsertParagraphAfter()
$scratchPara = $d.Paragraphs($d.Paragraphs.Count)

$scratchPara.Range.Style = "a6"
$balloonText = $d.Styles("a6")
$balloonText.NameLocal = "Balloon Text"
$balloonText.BaseStyle = "a"
$balloonText.Priority = 99
$balloonText.UnhideWhenUsed = $true
$balloonFont = $balloonText.Font
$balloonFont.Name = "Segoe UI"
$balloonFont.NameBi = "Segoe UI"
$balloonFont.Size = 9
$balloonFont.SizeBi = 9

$scratchPara.Range.Style = "a7"
$balloonTextChar = $d.Styles("a7")
$balloonTextChar.NameLocal = "Текст выноски Знак"
$balloonTextChar.BaseStyle = "a0"
$balloonTextChar.Priority = 99
$charFont = $balloonTextChar.Font
$charFont.Name = "Segoe UI"
$charFont.NameFarEast = "Times New Roman"
$charFont.NameBi = "Segoe UI"
$charFont.Size = 9
$charFont.SizeBi = 9
$charFont.LanguageIDFarEast = "ru-RU"

$balloonText.LinkStyle = "a7"
$balloonTextChar.LinkStyle = "a6"

$scratchPara.Range.Delete()
